$wb = $excel.ActiveWorkbook

# ---- Claims sheet: remove empty E/I placeholder cells (rows 51-77 E+I, rows 78-82 I only) ----
$claims = $wb.Worksheets.Item("Claims")
foreach ($r in 51..77) {
    $claims.Cells.Item($r, 5).ClearContents()
    $claims.Cells.Item($r, 9).ClearContents()
}
foreach ($r in 78..82) {
    $claims.Cells.Item($r, 9).ClearContents()
}

# ---- Claims sheet: add new rows 83-92 (reinstated v19 content) ----
# Row 83: C082
$claims.Cells.Item(83, 1).Value = "C082"
$claims.Cells.Item(83, 2).Value = "III.2"
$claims.Cells.Item(83, 3).Value = "The v19 distribution-gatekeeper layer indicates that corporatized veterinary networks, specialty retail chains, and scaled e-commerce platforms (including IVC Evidensia, Zooplus, PetSmart, and Musti) mediate demand access in the current market map [S116, Tab: Sheet1]."
$claims.Cells.Item(83, 4).Value = "S116"
$claims.Cells.Item(83, 5).Value = "III.2 paragraph 1"
$claims.Cells.Item(83, 6).Value = "N"
$claims.Cells.Item(83, 7).Value = "Y"
$claims.Cells.Item(83, 8).Value = "'2026-02-08"
$claims.Cells.Item(83, 10).Value = "[2026-02-08] Reinstated from v19 distribution-gatekeeper block."

# Row 84: C083
$claims.Cells.Item(84, 1).Value = "C083"
$claims.Cells.Item(84, 2).Value = "III.2"
$claims.Cells.Item(84, 3).Value = "The hypothesis that gatekeeper concentration reallocates margin toward owned-label portfolios remains unresolved and requires direct quantification [UNVERIFIED]."
$claims.Cells.Item(84, 4).Value = "UNVERIFIED"
$claims.Cells.Item(84, 5).Value = "III.2 paragraph 1"
$claims.Cells.Item(84, 6).Value = "N"
$claims.Cells.Item(84, 7).Value = "Y"
$claims.Cells.Item(84, 8).Value = "'2026-02-08"
$claims.Cells.Item(84, 10).Value = "[2026-02-08] Reinstated unresolved v19 hypothesis for follow-up validation."

# Row 85: C084
$claims.Cells.Item(85, 1).Value = "C084"
$claims.Cells.Item(85, 2).Value = "III.2"
$claims.Cells.Item(85, 3).Value = "Legacy v19 transaction history (Blue Buffalo, Neovia, Erber Group, Nom Nom, Aker BioMarine Feed Ingredients) is reinstated as an unresolved benchmark set pending direct source mapping [UNVERIFIED]."
$claims.Cells.Item(85, 4).Value = "UNVERIFIED"
$claims.Cells.Item(85, 5).Value = "III.2 paragraph 2; Table III.1"
$claims.Cells.Item(85, 6).Value = "N"
$claims.Cells.Item(85, 7).Value = "Y"
$claims.Cells.Item(85, 8).Value = "'2026-02-08"
$claims.Cells.Item(85, 10).Value = "[2026-02-08] Reinstated missing v19 transaction block with unresolved sourcing."

# Row 86: C085
$claims.Cells.Item(86, 1).Value = "C085"
$claims.Cells.Item(86, 2).Value = "III.2"
$claims.Cells.Item(86, 3).Value = "Legacy platform signals around Vetnique-Lintbells (YuMOVE) and FoodScience are mapped in the in-repo PE/VC portfolio source [S116, Tab: Sheet1]."
$claims.Cells.Item(86, 4).Value = "S116"
$claims.Cells.Item(86, 5).Value = "III.2 paragraph 2; Table III.1/Table III.2"
$claims.Cells.Item(86, 6).Value = "N"
$claims.Cells.Item(86, 7).Value = "Y"
$claims.Cells.Item(86, 8).Value = "'2026-02-08"
$claims.Cells.Item(86, 10).Value = "[2026-02-08] Added explicit source-backed sponsor mapping from internal portfolio workbook."

# Row 87: C086
$claims.Cells.Item(87, 1).Value = "C086"
$claims.Cells.Item(87, 2).Value = "III.2"
$claims.Cells.Item(87, 3).Value = "The v19 two-speed multiple framing (Pet roughly 15x-25x EBITDA vs Livestock/Feed roughly 8x-12x EBITDA) remains unresolved until directly source-bound [UNVERIFIED]."
$claims.Cells.Item(87, 4).Value = "UNVERIFIED"
$claims.Cells.Item(87, 5).Value = "III.2 paragraph 2"
$claims.Cells.Item(87, 6).Value = "N"
$claims.Cells.Item(87, 7).Value = "Y"
$claims.Cells.Item(87, 8).Value = "'2026-02-08"
$claims.Cells.Item(87, 10).Value = "[2026-02-08] Reinstated missing v19 valuation-band context; pending source-level validation."

# Row 88: C087
$claims.Cells.Item(88, 1).Value = "C087"
$claims.Cells.Item(88, 2).Value = "III.2"
$claims.Cells.Item(88, 3).Value = "The v19 investor-profile mapping is reinstated with fund-to-asset examples (JAB, Gryphon, MSCP, EQT, BC Partners, Cinven, Ani.VC) from the internal portfolio source [S116, Tab: Sheet1]."
$claims.Cells.Item(88, 4).Value = "S116"
$claims.Cells.Item(88, 5).Value = "III.2 paragraph 3; Table III.2"
$claims.Cells.Item(88, 6).Value = "N"
$claims.Cells.Item(88, 7).Value = "Y"
$claims.Cells.Item(88, 8).Value = "'2026-02-08"
$claims.Cells.Item(88, 10).Value = "[2026-02-08] Added back missing PE/VC investor-profile layer from v19."

# Row 89: C088
$claims.Cells.Item(89, 1).Value = "C088"
$claims.Cells.Item(89, 2).Value = "III.2"
$claims.Cells.Item(89, 3).Value = "Legacy buyer-mix assumptions (higher PE/financial participation in pet vs higher strategic participation in livestock/feed) remain unresolved pending deal-level coding [UNVERIFIED]."
$claims.Cells.Item(89, 4).Value = "UNVERIFIED"
$claims.Cells.Item(89, 5).Value = "III.2 paragraph 3"
$claims.Cells.Item(89, 6).Value = "N"
$claims.Cells.Item(89, 7).Value = "Y"
$claims.Cells.Item(89, 8).Value = "'2026-02-08"
$claims.Cells.Item(89, 10).Value = "[2026-02-08] Reinstated unresolved v19 buyer-mix split assumptions."

# Row 90: C089
$claims.Cells.Item(90, 1).Value = "C089"
$claims.Cells.Item(90, 2).Value = "III.2"
$claims.Cells.Item(90, 3).Value = "Legacy IPO optionality is retained, but explicit IPO comp tables are currently not source-mapped in this repository [UNVERIFIED]."
$claims.Cells.Item(90, 4).Value = "UNVERIFIED"
$claims.Cells.Item(90, 5).Value = "III.2 paragraph 3"
$claims.Cells.Item(90, 6).Value = "N"
$claims.Cells.Item(90, 7).Value = "Y"
$claims.Cells.Item(90, 8).Value = "'2026-02-08"
$claims.Cells.Item(90, 10).Value = "[2026-02-08] Reinstated missing IPO-language context with unresolved source status."

# Row 91: C090
$claims.Cells.Item(91, 1).Value = "C090"
$claims.Cells.Item(91, 2).Value = "III.2"
$claims.Cells.Item(91, 3).Value = "Legacy AUM tiers and investor-size rankings referenced in v19 remain unresolved until fund-level reporting sources are attached [UNVERIFIED]."
$claims.Cells.Item(91, 4).Value = "UNVERIFIED"
$claims.Cells.Item(91, 5).Value = "III.2 paragraph 3"
$claims.Cells.Item(91, 6).Value = "N"
$claims.Cells.Item(91, 7).Value = "Y"
$claims.Cells.Item(91, 8).Value = "'2026-02-08"
$claims.Cells.Item(91, 10).Value = "[2026-02-08] Reinstated unresolved v19 investor-size overlays."

# Row 92: C091
$claims.Cells.Item(92, 1).Value = "C091"
$claims.Cells.Item(92, 2).Value = "'01"
$claims.Cells.Item(92, 3).Value = "Global corporate, startup, and investor landscape visual added to executive summary, mapped from the internal PE/VC portfolio source and final composite image [S116, Tab: Sheet1; S121]."
$claims.Cells.Item(92, 4).Value = "S116, S121"
$claims.Cells.Item(92, 5).Value = "Figure ES-1"
$claims.Cells.Item(92, 6).Value = "N"
$claims.Cells.Item(92, 7).Value = "Y"
$claims.Cells.Item(92, 8).Value = "'2026-02-08"
$claims.Cells.Item(92, 10).Value = "[2026-02-08] Figure ES-1 insertion and provenance mapping."

# ---- Sources sheet: add new row 122 ----
$sources = $wb.Worksheets.Item("Sources")
$sources.Cells.Item(122, 1).Value = "S121"
$sources.Cells.Item(122, 2).Value = "Global Antigravity Landscape Final"
$sources.Cells.Item(122, 3).Value = "Internal Visual Composite"
$sources.Cells.Item(122, 4).Value = "_figures/exports/Global_Antigravity_Landscape_Final.png"
$sources.Cells.Item(122, 5).Value = "https://www.mapchart.net/world.html"
$sources.Cells.Item(122, 6).Value = "'2026-02-04"
$sources.Cells.Item(122, 7).Value = "'2026-02-08"
$sources.Cells.Item(122, 8).Value = "Codex"
$sources.Cells.Item(122, 9).Value = "Provenance traced in-repo: visual aligns with investor/company mapping from sources/internal/20260115_VC_PE_Portfolio.xlsx (Sheet1/Sheet2); base map tool watermark indicates mapchart.net. Lineage family found in _figures/exports: Global_Map_V10_VCPE.png -> Global_Map_V11_Final.png -> Global_Antigravity_Landscape*.png -> Global_Antigravity_Landscape_Final.png."

# ---- Figures sheet: add new row 46 ----
$figures = $wb.Worksheets.Item("Figures")
$figures.Cells.Item(46, 1).Value = "FIG-ES-1"
$figures.Cells.Item(46, 2).Value = "Global corporate, startup, and investor landscape by region."
$figures.Cells.Item(46, 3).Value = "Map"
$figures.Cells.Item(46, 4).Value = "S116, S121"
$figures.Cells.Item(46, 5).Value = "Figure 46"
$figures.Cells.Item(46, 6).Value = "Executive Summary"
$figures.Cells.Item(46, 7).Value = "Reviewed"
$figures.Cells.Item(46, 8).Value = "'2026-02-08"
$figures.Cells.Item(46, 9).Value = "Inserted in sections/01_executive_summary.md. Image file: figures/Global_Antigravity_Landscape_Final.png."
